# sp_AskBrent Check ID List - add the "Re-Compiles per Second" /
# "SQL Compilations/sec" checks (CheckID 25 & 26) and bump the title to v16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 29: CheckID 25, SQL Compilations/sec ---
$ws.Range("A29").Value = 25
$ws.Range("B29").Value = 250
$ws.Range("C29").Value = "Server Info"

# --- New row 30: CheckID 26, Re-Compiles per Second ---
$ws.Range("A30").Value = 26
$ws.Range("B30").Value = 250
$ws.Range("C30").Value = "Server Info"
$ws.Range("D30").Value = "Re-Compiles per Second"
$ws.Range("E30").Value = "http://BrentOzar.com/go/measure"
$ws.Hyperlinks.Add($ws.Range("E30"), "http://BrentOzar.com/go/measure")
$ws.Range("E30").Style = $ws.Range("E5").Style()

$ws.Range("D29").Value = "SQL Compilations/sec"
$ws.Range("E29").Value = "http://BrentOzar.com/go/measure"
$ws.Hyperlinks.Add($ws.Range("E29"), "http://BrentOzar.com/go/measure")
$ws.Range("E29").Style = $ws.Range("E5").Style()

# --- Update the title banner (row 1) last ---
$ws.Range("A1").Value = "sp_AskBrent Check ID List - v16 2015-07-18"
